# Automatische test-sync: 2025-06-29 15:20:50
#
# Adds the "Testmail #18" row to the Logs sheet, adds the corresponding
# "Planning / Afspraak" tally row to the Dashboard sheet, extends the
# conditional formatting ranges on the Logs sheet, and updates the
# category/value series references on the Dashboard bar chart.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 33
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A33").Value = "Ik wil een afspraak plannen, regel jij dat?"
$logs.Range("B33").Value = "mailmind.test@zohomail.eu"
$logs.Range("C33").Value = "Testmail #18: Ik wil een afspraak plannen, regel jij dat?"
$logs.Range("D33").Value = "Planning / Afspraak"

$answer33 = "Beste klant,`n" + `
    "Bedankt voor je e-mail. Helaas kan ik geen afspraken plannen via e-mail. Graag verzoek ik je om contact op te nemen met onze klantenservice of receptie om een afspraak te maken. Je kunt ons bereiken via [telefoonnummer] of [e-mailadres]. Onze medewerkers staan klaar om je verder te helpen.`n" + `
    "Met vriendelijke groet,`n" + `
    "[Naam] `n" + `
    "E-mailassistent bij [Bedrijfsnaam]"
$logs.Range("E33").Value = $answer33

$logs.Range("F33").Value = "2025-06-29 15:20:11"
$logs.Range("G33").Value = "Ja"
$logs.Range("H33").Value = "Ja"
$logs.Range("I33").Value = "Nee"

# ---------------------------------------------------------------------
# 2. Logs sheet: extend conditional formatting from row 32 to row 33
# ---------------------------------------------------------------------
$allCf = $logs.Cells.FormatConditions
$allCf.Item(1).ModifyAppliesToRange($logs.Range("D2:D33"))   # D2:D32 -> D2:D33 (group of 6 rules)
$allCf.Item(7).ModifyAppliesToRange($logs.Range("G2:G33"))   # G2:G32 -> G2:G33 (group of 2 rules)
$allCf.Item(9).ModifyAppliesToRange($logs.Range("H2:H33"))   # H2:H32 -> H2:H33
$allCf.Item(10).ModifyAppliesToRange($logs.Range("I2:I33"))  # I2:I32 -> I2:I33

# ---------------------------------------------------------------------
# 3. Dashboard sheet: append tally row 10
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("A10").Value = "Planning / Afspraak"
$dash.Range("B10").Value = 1

# ---------------------------------------------------------------------
# 4. Dashboard chart: extend category/value series ranges to row 10
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$series = $chartObj.Chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$10,Dashboard!`$B`$2:`$B`$10,1)"
